$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "D"
$ws.Range("E1").Value = "E"
$ws.Range("F1").Value = "F"
$ws.Range("G1").Value = "G"
$ws.Range("H1").Value = "H"
$ws.Range("I1").Value = "I"
$ws.Range("J1").Value = "J"
$ws.Range("K1").Value = "K"
